# Apply the "ok add zek items" edit to drama_zek_steal_bottle.
# 1. Change the flavor text of the resonance bottle (I8/J8): "void" -> "death".
# 2. Insert a setFlag + an extra modInvoke(complete_quest) row into the
#    "accept/refuse" branch (around old row 65), recording the player's
#    bottle_choice flag (0) before completing the refuse-path quest.
# 3. Insert a setFlag + an extra modInvoke(complete_quest) row into the
#    "ending" branch (around old row 90), recording the player's
#    bottle_choice flag (1) before completing the accept-path quest.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the narration text that describes the bottle -------------
$newText = "（リリィが作った『死の共鳴瓶』。彼女は満足げに、その器を受付の棚に飾っていた。）"
$ws.Cells.Item(8, 9).Value = $newText   # I8
$ws.Cells.Item(8, 10).Value = $newText  # J8

# --- 2. Insert the two new rows before the old row 65 --------------------
# Before: row 65 = modInvoke complete_quest(05_2_zek_steal_bottle_refuse), pc
#         row 66 = end
# After:  row 65 = setFlag chitsii.arena.player.bottle_choice,0
#         row 66 = modInvoke complete_quest(05_2_zek_steal_bottle), pc
#         row 67 = modInvoke complete_quest(05_2_zek_steal_bottle_refuse), pc  (shifted)
#         row 68 = end  (shifted)
$ws.Rows("65:66").Insert()

$ws.Cells.Item(65, 4).Value = "setFlag"
$ws.Cells.Item(65, 5).Value = "chitsii.arena.player.bottle_choice,0"

$ws.Cells.Item(66, 4).Value = "modInvoke"
$ws.Cells.Item(66, 5).Value = "complete_quest(05_2_zek_steal_bottle)"
$ws.Cells.Item(66, 6).Value = "pc"

# --- 3. Insert the two new rows before the old row 90 ---------------------
# After step 2, the old row 90 (modInvoke complete_quest(..._accept), pc)
# has shifted down by 2, to row 92; old row 91 (end) is now row 93.
# Before: row 92 = modInvoke complete_quest(05_2_zek_steal_bottle_accept), pc
#         row 93 = end
# After:  row 92 = setFlag chitsii.arena.player.bottle_choice,1
#         row 93 = modInvoke complete_quest(05_2_zek_steal_bottle), pc
#         row 94 = modInvoke complete_quest(05_2_zek_steal_bottle_accept), pc  (shifted)
#         row 95 = end  (shifted)
$ws.Rows("92:93").Insert()

$ws.Cells.Item(92, 4).Value = "setFlag"
$ws.Cells.Item(92, 5).Value = "chitsii.arena.player.bottle_choice,1"

$ws.Cells.Item(93, 4).Value = "modInvoke"
$ws.Cells.Item(93, 5).Value = "complete_quest(05_2_zek_steal_bottle)"
$ws.Cells.Item(93, 6).Value = "pc"
